# Fix NLI skipping some topics/activities (3/?)
#
# Appends 40 new evaluation rows (iCloud, iTunes, Refund topics) to the
# results sheet, right after the existing "Mac" rows (sheet previously
# ended at row 72). Also extends the dimension / conditional-formatting
# ranges to cover the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (A value, Label, Sentence, Optimal Threshold, MCC, Accuracy, Balanced Accuracy, F1, Items)
$newRows = @(
    @(71, "iCloud", "This example is icloud", 0.984, 0.1847239980152887, 0.955, 0.6489001692047377, 0.1818181818181818, 3),
    @(72, "iCloud", "This example is iCloud", 0.7999999999999999, 0.5071829918124481, 0.98, 0.8257191201353637, 0.5, 3),
    @(73, "iCloud", "This example is about icloud", 0.77, 0.6615905245346869, 0.99, 0.8307952622673435, 0.6666666666666666, 3),
    @(74, "iCloud", "This example is about iCloud", 0.93, 0.6615905245346869, 0.99, 0.8307952622673435, 0.6666666666666666, 3),
    @(75, "iCloud", "The example is icloud", 0.9979999999999999, 0.5744416851006197, 0.99, 0.6666666666666666, 0.5, 3),
    @(76, "iCloud", "The example is iCloud", 0.7899999999999999, 0.701702064128381, 0.985, 0.9923857868020305, 0.6666666666666666, 3),
    @(77, "iCloud", "The example is about icloud", 0.8099999999999999, 0.5700066912114538, 0.985, 0.8282571912013537, 0.5714285714285715, 3),
    @(78, "iCloud", "The example is about iCloud", 0.82, 0.7706546758708627, 0.99, 0.9949238578680203, 0.7499999999999999, 3),
    @(79, "iCloud", "The sentence is icloud", 0.999, 0, 0.985, 0.5, 0, 3),
    @(80, "iCloud", "The sentence is iCloud", 0.9099999999999998, 0.8144321109212623, 0.995, 0.8333333333333333, 0.8, 3),
    @(81, "iCloud", "The sentence is about icloud", 0.9399999999999999, 0.8144321109212623, 0.995, 0.8333333333333333, 0.8, 3),
    @(82, "iCloud", "The sentence is about iCloud", 0.9700000000000001, 0.6615905245346869, 0.99, 0.8307952622673435, 0.6666666666666666, 3),
    @(83, "iCloud", "The customer asks about icloud", 0.995, 0.8144321109212623, 0.995, 0.8333333333333333, 0.8, 3),
    @(84, "iCloud", "The customer asks about iCloud", 0.9899999999999999, 0.6615905245346869, 0.99, 0.8307952622673435, 0.6666666666666666, 3),
    @(85, "iTunes", "This example is itunes", 0.981, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(86, "iTunes", "This example is iTunes", 0.9099999999999999, 0.7448979591836735, 0.99, 0.8724489795918368, 0.75, 4),
    @(87, "iTunes", "This example is about itunes", 0.992, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(88, "iTunes", "This example is about iTunes", 0.989, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(89, "iTunes", "The example is itunes", 0.8899999999999999, 0.8123201004396181, 0.99, 0.9948979591836735, 0.8, 4),
    @(90, "iTunes", "The example is iTunes", 0.983, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(91, "iTunes", "The example is about itunes", 0.996, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(92, "iTunes", "The example is about iTunes", 0.994, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(93, "iTunes", "The sentence is itunes", 0.87, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(94, "iTunes", "The sentence is iTunes", 0.86, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(95, "iTunes", "The sentence is about itunes", 0.993, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(96, "iTunes", "The sentence is about iTunes", 0.994, 0.8638245732792135, 0.995, 0.875, 0.8571428571428571, 4),
    @(97, "iTunes", "The customer asks about itunes", 0.997, 0.7448979591836735, 0.99, 0.8724489795918368, 0.75, 4),
    @(98, "iTunes", "The customer asks about iTunes", 0.995, 0.7448979591836735, 0.99, 0.8724489795918368, 0.75, 4),
    @(99, "Refund", "This example is refund", 0.8899999999999999, 0.7706546758708627, 0.99, 0.8, 0.7499999999999999, 5),
    @(100, "Refund", "This example is about refund", 0.982, 1, 1, 1, 1, 5),
    @(101, "Refund", "This example is about a refund", 0.984, 1, 1, 1, 1, 5),
    @(102, "Refund", "The example is refund", 0.8899999999999999, 0.910527225875113, 0.995, 0.9974358974358974, 0.9090909090909091, 5),
    @(103, "Refund", "The example is about refund", 0.99, 0.910527225875113, 0.995, 0.9974358974358974, 0.9090909090909091, 5),
    @(104, "Refund", "The example is about a refund", 0.991, 0.910527225875113, 0.995, 0.9974358974358974, 0.9090909090909091, 5),
    @(105, "Refund", "The sentence is refund", 0.9099999999999999, 0.7706546758708627, 0.99, 0.8, 0.7499999999999999, 5),
    @(106, "Refund", "The sentence is about refund", 0.9890000000000001, 0.910527225875113, 0.995, 0.9974358974358974, 0.9090909090909091, 5),
    @(107, "Refund", "The sentence is about a refund", 0.9949999999999999, 1, 1, 1, 1, 5),
    @(108, "Refund", "The customer asks about refund", 0.998, 1, 1, 1, 1, 5),
    @(109, "Refund", "The customer asks about a refund", 0.998, 1, 1, 1, 1, 5),
    @(110, "Refund", "The customer wants a refund", 0.9960000000000001, 1, 1, 1, 1, 5)
)

$startRow = 73

# Give the new A-column cells the same style (bold, centered, bordered) as
# the existing index column before filling in the values.
$ws.Range("A2").Copy()
$endRow = $startRow + $newRows.Count - 1
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

# Extend the color-scale conditional formatting (originally E2:E72 /
# F2:F72 / G2:G72 / H2:H72) to cover the newly added rows, preserving the
# existing rule (priority/colors/stops) rather than recreating it.
$lastRow = $endRow
$ws.Range("E2:E72").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E2:E$lastRow"))
$ws.Range("F2:F72").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("F2:F$lastRow"))
$ws.Range("G2:G72").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G$lastRow"))
$ws.Range("H2:H72").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H$lastRow"))
